$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 6089.421
$ws.Range("I113").Value = 4892.385
$ws.Range("J113").Value = 8683
$ws.Range("K113").Value = 4892.385
$ws.Range("L113").Value = 8683
$ws.Range("M113").Value = -1638.385
$ws.Range("N113").Value = -15191
$ws.Range("H132").Value = 2606.6616
$ws.Range("I132").Value = 2350.7705
$ws.Range("J132").Value = 6509
$ws.Range("K132").Value = 7052.3115
$ws.Range("L132").Value = 19527
$ws.Range("M132").Value = -4522.3115
$ws.Range("N132").Value = -24587
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2071.8293
$ws.Range("I2").Value = 2009
$ws.Range("K2").Value = 2009
$ws.Range("M2").Value = -1896
$ws.Range("H32").Value = 6307.465
$ws.Range("I32").Value = 6159.2646
$ws.Range("K32").Value = 6159.2646
$ws.Range("M32").Value = -5872.2646
$ws.Range("H61").Value = 9492.440000000001
$ws.Range("I61").Value = 10859.368
$ws.Range("K61").Value = 10859.368
$ws.Range("M61").Value = -10647.368
$ws.Range("H74").Value = 6120.0835
$ws.Range("I74").Value = 10508.1
$ws.Range("J74").Value = 2985.7856
$ws.Range("K74").Value = 10508.1
$ws.Range("L74").Value = 2985.7856
$ws.Range("M74").Value = -9634.1
$ws.Range("N74").Value = -4733.7856
$ws.Range("H77").Value = 6120.0835
$ws.Range("I77").Value = 10508.1
$ws.Range("J77").Value = 2985.7856
$ws.Range("K77").Value = 52540.5
$ws.Range("L77").Value = 14928.928
$ws.Range("M77").Value = -48172.5
$ws.Range("N77").Value = -23664.928
$ws.Range("H116").Value = 2071.8293
$ws.Range("I116").Value = 2009
$ws.Range("K116").Value = 2009
$ws.Range("M116").Value = 285
$ws.Range("H132").Value = 2941
$ws.Range("I132").Value = 2235.1538
$ws.Range("K132").Value = 6705.4614
$ws.Range("M132").Value = -4175.4614
$ws.Range("H136").Value = 9492.440000000001
$ws.Range("I136").Value = 10859.368
$ws.Range("K136").Value = 32578.104
$ws.Range("M136").Value = -30028.104
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2071.8293
$ws.Range("I3").Value = 2009
$ws.Range("K3").Value = 2009
$ws.Range("M3").Value = -1895
$ws.Range("H94").Value = 1944.5428
$ws.Range("I94").Value = 1011.96
$ws.Range("K94").Value = 1011.96
$ws.Range("M94").Value = -560.96
$ws.Range("H95").Value = 35249.5
$ws.Range("J95").Value = 35249.5
$ws.Range("L95").Value = 35249.5
$ws.Range("N95").Value = -40741.5
$ws.Range("H99").Value = 8936.166999999999
$ws.Range("I99").Value = 10910.423
$ws.Range("J99").Value = 3803.1
$ws.Range("K99").Value = 10910.423
$ws.Range("L99").Value = 3803.1
$ws.Range("M99").Value = -9412.423000000001
$ws.Range("N99").Value = -6799.1
$ws.Range("H107").Value = 1684.2307
$ws.Range("I107").Value = 1743.0435
$ws.Range("K107").Value = 1743.0435
$ws.Range("M107").Value = 176.9565
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4565.3193
$ws.Range("I31").Value = 4475.5264
$ws.Range("J31").Value = 4944.4443
$ws.Range("K31").Value = 4475.5264
$ws.Range("L31").Value = 4944.4443
$ws.Range("M31").Value = -4180.5264
$ws.Range("N31").Value = -5534.4443
$ws.Range("H34").Value = 4565.3193
$ws.Range("I34").Value = 4475.5264
$ws.Range("J34").Value = 4944.4443
$ws.Range("K34").Value = 4475.5264
$ws.Range("L34").Value = 4944.4443
$ws.Range("M34").Value = -4273.5264
$ws.Range("N34").Value = -5348.4443
$ws.Range("H58").Value = 3588.8
$ws.Range("I58").Value = 4086.2307
$ws.Range("K58").Value = 4086.2307
$ws.Range("M58").Value = -3883.2307
$ws.Range("H105").Value = 7838.933
$ws.Range("I105").Value = 11260.8
$ws.Range("K105").Value = 11260.8
$ws.Range("M105").Value = -9513.799999999999
$ws.Range("H117").Value = 22322.8
$ws.Range("J117").Value = 21653.5
$ws.Range("L117").Value = 21653.5
$ws.Range("N117").Value = -30831.5
$ws.Range("H132").Value = 2037.6875
$ws.Range("I132").Value = 1614.5
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 4843.5
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -2313.5
$ws.Range("N132").Value = -20060
$ws.Range("H134").Value = 10091.071
$ws.Range("I134").Value = 16659.375
$ws.Range("K134").Value = 49978.125
$ws.Range("M134").Value = -47443.125
$ws.Range("H136").Value = 3588.8
$ws.Range("I136").Value = 4086.2307
$ws.Range("K136").Value = 12258.6921
$ws.Range("M136").Value = -9708.6921
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 25093.35
$ws.Range("I11").Value = 29498.941
$ws.Range("J11").Value = 128.33333
$ws.Range("K11").Value = 88496.823
$ws.Range("L11").Value = 384.99999
$ws.Range("M11").Value = -88356.823
$ws.Range("N11").Value = -664.99999
$ws.Range("H22").Value = 876.1
$ws.Range("I22").Value = 896.2857
$ws.Range("J22").Value = 829
$ws.Range("K22").Value = 2688.8571
$ws.Range("L22").Value = 2487
$ws.Range("M22").Value = -2519.8571
$ws.Range("N22").Value = -2825
$ws.Range("H27").Value = 876.1
$ws.Range("I27").Value = 896.2857
$ws.Range("J27").Value = 829
$ws.Range("K27").Value = 2688.8571
$ws.Range("L27").Value = 2487
$ws.Range("M27").Value = -2586.8571
$ws.Range("N27").Value = -2691
$ws.Range("H96").Value = 3995
$ws.Range("J96").Value = 3995
$ws.Range("L96").Value = 11985
$ws.Range("N96").Value = -16103
$ws.Range("H98").Value = 2199.3333
$ws.Range("J98").Value = 1971.1428
$ws.Range("L98").Value = 5913.428400000001
$ws.Range("N98").Value = -8909.428400000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8812.954
$ws.Range("I70").Value = 7100.3335
$ws.Range("K70").Value = 7100.3335
$ws.Range("M70").Value = -6830.3335
$ws.Range("H73").Value = 8812.954
$ws.Range("I73").Value = 7100.3335
$ws.Range("K73").Value = 7100.3335
$ws.Range("M73").Value = -6164.3335
$ws.Range("H113").Value = 9486.933999999999
$ws.Range("I113").Value = 14538.625
$ws.Range("K113").Value = 14538.625
$ws.Range("M113").Value = -12368.625
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H126").Value = 14278.818
$ws.Range("I126").Value = 41994.332
$ws.Range("K126").Value = 125982.996
$ws.Range("M126").Value = -123512.996
$ws.Range("H132").Value = 6531.4707
$ws.Range("I132").Value = 6735.6665
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 20206.9995
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -17676.9995
$ws.Range("N132").Value = -20060
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H132").Value = 1492874.2
$ws.Range("I132").Value = 2485058.8
$ws.Range("J132").Value = 4597.5
$ws.Range("K132").Value = 7455176.399999999
$ws.Range("L132").Value = 13792.5
$ws.Range("M132").Value = -7452646.399999999
$ws.Range("N132").Value = -18852.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 673
$ws.Range("I4").Value = 110
$ws.Range("J4").Value = 1024.875
$ws.Range("K4").Value = 110
$ws.Range("L4").Value = 1024.875
$ws.Range("M4").Value = 3
$ws.Range("N4").Value = -1250.875
$ws.Range("H42").Value = 25000
$ws.Range("J42").Value = 25000
$ws.Range("L42").Value = 25000
$ws.Range("N42").Value = -25756
$ws.Range("H112").Value = 49966.332
$ws.Range("J112").Value = 49966.332
$ws.Range("L112").Value = 49966.332
$ws.Range("N112").Value = -52920.332
$ws.Range("H135").Value = 51500
$ws.Range("J135").Value = 51500
$ws.Range("L135").Value = 51500
$ws.Range("N135").Value = -61640
$ws.Range("H136").Value = 446803.66
$ws.Range("I136").Value = 597109.4
$ws.Range("J136").Value = 12587.111
$ws.Range("K136").Value = 1791328.2
$ws.Range("L136").Value = 37761.333
$ws.Range("M136").Value = -1788778.2
$ws.Range("N136").Value = -42861.333
